$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update intro text (A2) with new "as at" date
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 29 August 2025"

# Update data table rows 5-61
$ws.Range("A5").Value = "01 Sep 2025"
$ws.Range("B5").Value = "Civil justice statistics: April to June 2025"
$ws.Range("C5").Value = "4 September 2025"
$ws.Range("D5").Value = "provisional"
$ws.Range("E5").Value = 36
$ws.Range("F5").Value = "standard"

$ws.Range("A6").Value = "08 Sep 2025"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = 37
$ws.Range("F6").Value = ""

$ws.Range("A7").Value = "15 Sep 2025"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = 38
$ws.Range("F7").Value = ""

$ws.Range("A8").Value = "22 Sep 2025"
$ws.Range("B8").Value = "Family court statistics quarterly: April to June 2025"
$ws.Range("C8").Value = "25 September 2025"
$ws.Range("D8").Value = "provisional"
$ws.Range("E8").Value = 39
$ws.Range("F8").Value = "standard"

$ws.Range("A9").Value = "22 Sep 2025"
$ws.Range("B9").Value = "Education and Accredited Programme Statistics 2024 to 2025"
$ws.Range("C9").Value = "25 September 2025"
$ws.Range("D9").Value = "provisional"
$ws.Range("E9").Value = 39
$ws.Range("F9").Value = "standard"

$ws.Range("A10").Value = "22 Sep 2025"
$ws.Range("B10").Value = "Legal aid statistics quarterly: April to June 2025"
$ws.Range("C10").Value = "25 September 2025"
$ws.Range("D10").Value = "provisional"
$ws.Range("E10").Value = 39
$ws.Range("F10").Value = "standard"

$ws.Range("A11").Value = "29 Sep 2025"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = 40
$ws.Range("F11").Value = ""

$ws.Range("A12").Value = "06 Oct 2025"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 41
$ws.Range("F12").Value = ""

$ws.Range("A13").Value = "13 Oct 2025"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 42
$ws.Range("F13").Value = ""

$ws.Range("A14").Value = "20 Oct 2025"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = 43
$ws.Range("F14").Value = ""

$ws.Range("A15").Value = "27 Oct 2025"
$ws.Range("B15").Value = "Justice data lab statistics: October 2025"
$ws.Range("C15").Value = "30 October 2025"
$ws.Range("D15").Value = "provisional"
$ws.Range("E15").Value = 44
$ws.Range("F15").Value = "standard"

$ws.Range("A16").Value = "27 Oct 2025"
$ws.Range("B16").Value = "Safety in custody: quarterly update to June 2025"
$ws.Range("C16").Value = "30 October 2025"
$ws.Range("D16").Value = "provisional"
$ws.Range("E16").Value = 44
$ws.Range("F16").Value = "standard"

$ws.Range("A17").Value = "27 Oct 2025"
$ws.Range("B17").Value = "Justice data lab statistics: October 2025"
$ws.Range("C17").Value = "30 October 2025"
$ws.Range("D17").Value = "provisional"
$ws.Range("E17").Value = 44
$ws.Range("F17").Value = "standard"

$ws.Range("A18").Value = "27 Oct 2025"
$ws.Range("B18").Value = "Proven reoffending statistics: October to December 2023"
$ws.Range("C18").Value = "30 October 2025"
$ws.Range("D18").Value = "provisional"
$ws.Range("E18").Value = 44
$ws.Range("F18").Value = "standard"

$ws.Range("A19").Value = "27 Oct 2025"
$ws.Range("B19").Value = "Deaths of offenders supervised in the community, England and Wales, 2024/2025"
$ws.Range("C19").Value = "30 October 2025"
$ws.Range("D19").Value = "provisional"
$ws.Range("E19").Value = 44
$ws.Range("F19").Value = "standard"

$ws.Range("A20").Value = "03 Nov 2025"
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = 45
$ws.Range("F20").Value = ""

$ws.Range("A21").Value = "10 Nov 2025"
$ws.Range("B21").Value = "Prison Population Projections: 2025 to 2030"
$ws.Range("C21").Value = "13 November 2025"
$ws.Range("D21").Value = "provisional"
$ws.Range("E21").Value = 46
$ws.Range("F21").Value = "standard"

$ws.Range("A22").Value = "17 Nov 2025"
$ws.Range("B22").Value = "Knife and Offensive Weapon Sentencing Statistics:  April to June 2025"
$ws.Range("C22").Value = "20 November 2025"
$ws.Range("D22").Value = "provisional"
$ws.Range("E22").Value = 47
$ws.Range("F22").Value = "standard"

$ws.Range("A23").Value = "17 Nov 2025"
$ws.Range("B23").Value = " HM Prison and Probation Service workforce quarterly: September 2025"
$ws.Range("C23").Value = "20 November 2025"
$ws.Range("D23").Value = "provisional"
$ws.Range("E23").Value = 47
$ws.Range("F23").Value = "standard"

$ws.Range("A24").Value = "24 Nov 2025"
$ws.Range("B24").Value = "Her Majesty’s Prison and Probation Service offender equalities report: 2024 to 2025"
$ws.Range("C24").Value = "27 November 2025"
$ws.Range("D24").Value = "provisional"
$ws.Range("E24").Value = 48
$ws.Range("F24").Value = "standard"

$ws.Range("A25").Value = "24 Nov 2025"
$ws.Range("B25").Value = "Ethnicity and the Criminal Justice System 2024"
$ws.Range("C25").Value = "27 November 2025"
$ws.Range("D25").Value = "provisional"
$ws.Range("E25").Value = 48
$ws.Range("F25").Value = "standard"

$ws.Range("A26").Value = "01 Dec 2025"
$ws.Range("B26").Value = " Civil justice statistics: July to September 2025"
$ws.Range("C26").Value = "4 December 2025"
$ws.Range("D26").Value = "provisional"
$ws.Range("E26").Value = 49
$ws.Range("F26").Value = "standard"

$ws.Range("A27").Value = "08 Dec 2025"
$ws.Range("B27").Value = "Tribunals statistics quarterly: April to September 2025"
$ws.Range("C27").Value = "11 December 2025"
$ws.Range("D27").Value = "provisional"
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = "standard"

$ws.Range("A28").Value = "15 Dec 2025"
$ws.Range("B28").Value = "Family court statistics quarterly: July to September 2025"
$ws.Range("C28").Value = "18 December 2025"
$ws.Range("D28").Value = "provisional"
$ws.Range("E28").Value = 51
$ws.Range("F28").Value = "standard"

$ws.Range("A29").Value = "15 Dec 2025"
$ws.Range("B29").Value = "Criminal court statistics quarterly: July to September 2025"
$ws.Range("C29").Value = "18 December 2025"
$ws.Range("D29").Value = "provisional"
$ws.Range("E29").Value = 51
$ws.Range("F29").Value = "standard"

$ws.Range("A30").Value = "15 Dec 2025"
$ws.Range("B30").Value = "Legal aid statistics quarterly: July to September 2025"
$ws.Range("C30").Value = "18 December 2025"
$ws.Range("D30").Value = "provisional"
$ws.Range("E30").Value = 51
$ws.Range("F30").Value = "standard"

$ws.Range("A31").Value = "22 Dec 2025"
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = 52
$ws.Range("F31").Value = ""

$ws.Range("A32").Value = "05 Jan 2026"
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = ""
$ws.Range("E32").Value = 2
$ws.Range("F32").Value = ""

$ws.Range("A33").Value = "12 Jan 2026"
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = 3
$ws.Range("F33").Value = ""

$ws.Range("A34").Value = "19 Jan 2026"
$ws.Range("B34").Value = ""
$ws.Range("C34").Value = ""
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = ""

$ws.Range("A35").Value = "26 Jan 2026"
$ws.Range("B35").Value = ""
$ws.Range("C35").Value = ""
$ws.Range("D35").Value = ""
$ws.Range("E35").Value = 5
$ws.Range("F35").Value = ""

$ws.Range("A36").Value = "02 Feb 2026"
$ws.Range("B36").Value = ""
$ws.Range("C36").Value = ""
$ws.Range("D36").Value = ""
$ws.Range("E36").Value = 6
$ws.Range("F36").Value = ""

$ws.Range("A37").Value = "09 Feb 2026"
$ws.Range("B37").Value = ""
$ws.Range("C37").Value = ""
$ws.Range("D37").Value = ""
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = ""

$ws.Range("A38").Value = "16 Feb 2026"
$ws.Range("B38").Value = ""
$ws.Range("C38").Value = ""
$ws.Range("D38").Value = ""
$ws.Range("E38").Value = 8
$ws.Range("F38").Value = ""

$ws.Range("A39").Value = "23 Feb 2026"
$ws.Range("B39").Value = ""
$ws.Range("C39").Value = ""
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = ""

$ws.Range("A40").Value = "02 Mar 2026"
$ws.Range("B40").Value = ""
$ws.Range("C40").Value = ""
$ws.Range("D40").Value = ""
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = ""

$ws.Range("A41").Value = "09 Mar 2026"
$ws.Range("B41").Value = ""
$ws.Range("C41").Value = ""
$ws.Range("D41").Value = ""
$ws.Range("E41").Value = 11
$ws.Range("F41").Value = ""

$ws.Range("A42").Value = "16 Mar 2026"
$ws.Range("B42").Value = ""
$ws.Range("C42").Value = ""
$ws.Range("D42").Value = ""
$ws.Range("E42").Value = 12
$ws.Range("F42").Value = ""

$ws.Range("A43").Value = "23 Mar 2026"
$ws.Range("B43").Value = ""
$ws.Range("C43").Value = ""
$ws.Range("D43").Value = ""
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = ""

$ws.Range("A44").Value = "30 Mar 2026"
$ws.Range("B44").Value = ""
$ws.Range("C44").Value = ""
$ws.Range("D44").Value = ""
$ws.Range("E44").Value = 14
$ws.Range("F44").Value = ""

$ws.Range("A45").Value = "06 Apr 2026"
$ws.Range("B45").Value = ""
$ws.Range("C45").Value = ""
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = ""

$ws.Range("A46").Value = "13 Apr 2026"
$ws.Range("B46").Value = ""
$ws.Range("C46").Value = ""
$ws.Range("D46").Value = ""
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = ""

$ws.Range("A47").Value = "20 Apr 2026"
$ws.Range("B47").Value = ""
$ws.Range("C47").Value = ""
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = 17
$ws.Range("F47").Value = ""

$ws.Range("A48").Value = "27 Apr 2026"
$ws.Range("B48").Value = ""
$ws.Range("C48").Value = ""
$ws.Range("D48").Value = ""
$ws.Range("E48").Value = 18
$ws.Range("F48").Value = ""

$ws.Range("A49").Value = "04 May 2026"
$ws.Range("B49").Value = ""
$ws.Range("C49").Value = ""
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = 19
$ws.Range("F49").Value = ""

$ws.Range("A50").Value = "11 May 2026"
$ws.Range("B50").Value = ""
$ws.Range("C50").Value = ""
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = 20
$ws.Range("F50").Value = ""

$ws.Range("A51").Value = "18 May 2026"
$ws.Range("B51").Value = ""
$ws.Range("C51").Value = ""
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = 21
$ws.Range("F51").Value = ""

$ws.Range("A52").Value = "25 May 2026"
$ws.Range("B52").Value = ""
$ws.Range("C52").Value = ""
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = 22
$ws.Range("F52").Value = ""

$ws.Range("A53").Value = "01 Jun 2026"
$ws.Range("B53").Value = ""
$ws.Range("C53").Value = ""
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = 23
$ws.Range("F53").Value = ""

$ws.Range("A54").Value = "08 Jun 2026"
$ws.Range("B54").Value = ""
$ws.Range("C54").Value = ""
$ws.Range("D54").Value = ""
$ws.Range("E54").Value = 24
$ws.Range("F54").Value = ""

$ws.Range("A55").Value = "15 Jun 2026"
$ws.Range("B55").Value = ""
$ws.Range("C55").Value = ""
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = 25
$ws.Range("F55").Value = ""

$ws.Range("A56").Value = "22 Jun 2026"
$ws.Range("B56").Value = ""
$ws.Range("C56").Value = ""
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = 26
$ws.Range("F56").Value = ""

$ws.Range("A57").Value = "29 Jun 2026"
$ws.Range("B57").Value = ""
$ws.Range("C57").Value = ""
$ws.Range("D57").Value = ""
$ws.Range("E57").Value = 27
$ws.Range("F57").Value = ""

$ws.Range("A58").Value = "06 Jul 2026"
$ws.Range("B58").Value = ""
$ws.Range("C58").Value = ""
$ws.Range("D58").Value = ""
$ws.Range("E58").Value = 28
$ws.Range("F58").Value = ""

$ws.Range("A59").Value = "13 Jul 2026"
$ws.Range("B59").Value = ""
$ws.Range("C59").Value = ""
$ws.Range("D59").Value = ""
$ws.Range("E59").Value = 29
$ws.Range("F59").Value = ""

$ws.Range("A60").Value = "20 Jul 2026"
$ws.Range("B60").Value = ""
$ws.Range("C60").Value = ""
$ws.Range("D60").Value = ""
$ws.Range("E60").Value = 30
$ws.Range("F60").Value = ""

$ws.Range("A61").Value = "27 Jul 2026"
$ws.Range("B61").Value = "Community Performance Annual, update to March 2026"
$ws.Range("C61").Value = "30 July 2026"
$ws.Range("D61").Value = "confirmed"
$ws.Range("E61").Value = 31
$ws.Range("F61").Value = "standard"

# Remove now-unused trailing rows 62-64
$ws.Range("A62:F64").EntireRow.Delete()

# Fix conditional formatting ranges to match the new table extent
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo.Address()
    if ($addr -eq "`$A`$5:`$F`$64") {
        $fc.ModifyAppliesToRange($ws.Range("A5:F61"))
    } elseif ($addr -eq "`$A`$5:`$A`$64") {
        $fc.ModifyAppliesToRange($ws.Range("A5:A61"))
    }
}

Write-Output "Edit complete"